# Processing excel data and display view added
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Fill the new "Code" column (F) top to bottom
$ws.Range("F1").Value = "Code"
$ws.Range("F2").Value = "'001"
$ws.Range("F3").Value = "'002"
$ws.Range("F4").Value = "'003"

# Fill the new "Phone" column (G) top to bottom
$ws.Range("G1").Value = "Phone"
$ws.Range("G2").Value = "'08031812695"
$ws.Range("G3").Value = "'08069784914"

# Fill the second "Phone" column (H) top to bottom
$ws.Range("H1").Value = "Phone"
$ws.Range("H2").Value = "'+2348031812689"
$ws.Range("H4").Value = "'07029478943"

# Activate Sheet2 and move the selection to H3 (matches the saved view state)
$ws.Activate() | Out-Null
$ws.Range("H3").Select() | Out-Null
